$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.995.29'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.925.20'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '590.58'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').Value = '147.44'
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.507'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '33.72'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '3.409.13'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '60.943.59'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.925.16'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '6.70'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').Value = '432.50'
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').Value = '13.43'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '81.41'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').Value = '10.90'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '11.91'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +5.35%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '7.02'
$ws.Range('E30').Value = '  -2.81%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '0.110'
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  -5.12%  '
$ws.Range('D40').Value = '8.56'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').Value = '41.44'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('D42').Value = '0.283'
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').Value = '378.06'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0345'
$ws.Range('E44').Value = '  -1.09%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.708.08'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('D46').Value = '134.18'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '23.91'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('E49').Value = '  -0.59%  '
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('E51').Value = '  -0.61%  '
